# Append 6 new rows (15-20) to the "Captured_Values" sheet, following the
# existing pattern already present in rows 2-14: column A holds the number
# 123456789 and column B holds the text "Real Programmers Count 0123456789
# From Zero".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Captured_Values")

$value = "Real Programmers Count 0123456789 From Zero"

for ($row = 15; $row -le 20; $row++) {
    $ws.Cells.Item($row, 1).Value = 123456789
    $ws.Cells.Item($row, 2).Value = $value
}
